{"js": "// Remove the auto-managed \"_GoBack\" bookmark up front so our OOXML block\n// (which re-creates it in its new location) doesn't collide with the one\n// still sitting in the trailing paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nconst body = context.document.body;\n\n// Locate the \"docker-compose down \" paragraph via search (robust to any\n// paragraph-index drift) and take the body's very last paragraph (the\n// paragraph that only holds the _GoBack bookmark) as the end of the range\n// we are going to replace.\nconst hits = body.search(\"docker-compose down \", { matchCase: true });\nhits.load(\"items\");\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find 'docker-compose down ' paragraph\");\n}\n\nconst startPara = hits.items[0].paragraphs.getFirst();\nconst lastPara = paragraphs.items[paragraphs.items.length - 1];\n\nconst startRange = startPara.getRange(Word.RangeLocation.start);\nconst endRange = lastPara.getRange(Word.RangeLocation.end);\nconst targetRange = startRange.expandTo(endRange);\n\n// Flat-OPC wrapped OOXML for the full replacement block: the new\n// \"PING .../Install ping/apt-get ...\" paragraphs followed by the\n// (unchanged) \"docker-compose down \" / \"**must be in same directory...\"\n// paragraphs and a trailing empty paragraph, mirroring the target diff.\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\" pkg:padding=\"512\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:lastRenderedPageBreak/><w:t xml:space=\"preserve\">PING </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>webapi</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> from </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>webapp</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n          </w:p>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">Install ping in container: </w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr>\n            <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>apt-get update</w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr>\n            <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space=\"preserve\">apt-get install </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>iputils</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>-ping</w:t></w:r>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n          </w:p>\n          <w:p/>\n          <w:p>\n            <w:pPr><w:rPr><w:b/></w:rPr></w:pPr>\n            <w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">docker-compose down </w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:r><w:t>**must be in same directory as compose file</w:t></w:r>\n          </w:p>\n          <w:p/>\n          <w:sectPr/>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ntargetRange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the auto-managed \"_GoBack\" bookmark up front so the copy we\n# re-create in its new location (inside the OOXML block below) doesn't\n# collide with the stale one still anchored on the trailing paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Locate the \"docker-compose down \" paragraph via Find (robust to any\n# paragraph-index drift) and the paragraph just before the very last\n# (now bookmark-less) paragraph -- i.e. \"**must be in same directory\n# as compose file\" -- to build the range we are going to replace.\n$find = $d.Content\n$find.Find.Execute(\"docker-compose down \")\n$startOffset = $find.Start\n\n$countBefore = $d.Paragraphs.Count\n$secondToLastPara = $d.Paragraphs.Item($countBefore - 1)\n$endOffset = $secondToLastPara.Range.End\n\n$targetRange = $d.Range($startOffset, $endOffset)\n\n# Flat-OPC wrapped OOXML for the replacement block: the new\n# \"PING .../Install ping/apt-get ...\" paragraphs followed by the\n# (unchanged) \"docker-compose down \" / \"**must be in same directory...\"\n# paragraphs, mirroring the target diff. The trailing empty paragraph\n# that used to hold the bookmark is left alone (it still exists right\n# after $endOffset, untouched by this replace).\n$ooxml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\" pkg:padding=\"512\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:lastRenderedPageBreak/><w:t xml:space=\"preserve\">PING </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>webapi</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> from </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>webapp</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n          </w:p>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">Install ping in container: </w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr>\n            <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>apt-get update</w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr>\n            <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space=\"preserve\">apt-get install </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>iputils</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>-ping</w:t></w:r>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n          </w:p>\n          <w:p/>\n          <w:p>\n            <w:pPr><w:rPr><w:b/></w:rPr></w:pPr>\n            <w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">docker-compose down </w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:r><w:t>**must be in same directory as compose file</w:t></w:r>\n          </w:p>\n          <w:sectPr/>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$targetRange.InsertXML($ooxml)\n"}
